# Fix background and migrate to new template
# Update timing statistics (Tiempo_Minimo, Tiempo_Maximo, Tiempo_Promedio)
# for rows 2 and 3 on the "Data" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Highs-BigM (100,100))
$ws.Range("E2").Value = 0.000285655
$ws.Range("F2").Value = 0.02652873
$ws.Range("G2").Value = 0.0004965854824447335

# Row 3 (NO_SOLUTION)
$ws.Range("E3").Value = 0.004229489
$ws.Range("F3").Value = 0.010075473
$ws.Range("G3").Value = 0.00501175204321608
